# Insert one new weekly record for "Femacal de La Calera" / "Achicoria".
# This pushes the existing rows 177:228 down to 178:229 and fills the
# freshly opened row 177 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177 (shifts old 177..228 down to 178..229, carrying
# their original values/formatting with them automatically).
$ws.Rows("177:177").Insert()

# Capture the date column's number format from a neighboring (already
# shifted) row so the new row keeps the same "YYYY-MM-DD HH:MM:SS" style.
$dateFormat = $ws.Range("D178").NumberFormat

# Populate the new row 177 with the new week's values.
$ws.Range("A177").Value2 = 3
$ws.Range("B177").Value2 = "Femacal de La Calera"
$ws.Range("C177").Value2 = "Coquimbo"
$ws.Range("D177").Value2 = 44841
$ws.Range("D177").NumberFormat = $dateFormat
$ws.Range("E177").Value2 = 5
$ws.Range("F177").Value2 = 100112010
$ws.Range("G177").Value2 = "Achicoria"
$ws.Range("H177").Value2 = "Sin especificar"
$ws.Range("I177").Value2 = "Primera"
$ws.Range("J177").Value2 = 65
$ws.Range("K177").Value2 = 6000
$ws.Range("L177").Value2 = 6000
$ws.Range("M177").Value2 = 6000
$ws.Range("N177").Value2 = "$/caja 16 unidades"
$ws.Range("O177").Value2 = "Provincia de Quillota"
$ws.Range("P177").Value2 = 375
$ws.Range("Q177").Value2 = 16
$ws.Range("R177").Value2 = "Hortaliza"
